$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "quality_comparison"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("quality_comparison")

# C1 -> border: top + bottom only (no left/right)
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item(8).LineStyle = 1
$c1.Borders.Item(9).LineStyle = 1

# D1 -> border: top + bottom + right (no left)
$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.Item(8).LineStyle = 1
$d1.Borders.Item(10).LineStyle = 1
$d1.Borders.Item(9).LineStyle = 1

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# ---------------------------------------------------------------------------
# Sheet "computational_comparison"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("computational_comparison")

# Re-use the formatting derived above (copy/paste-formats) instead of
# re-deriving it via Borders again, which would leave stray/duplicate
# style entries behind in the shared style table.
$c1.Copy()
$c1b = $ws2.Range("C1")
$c1b.PasteSpecial(-4122)

$d1.Copy()
$d1b = $ws2.Range("D1")
$d1b.PasteSpecial(-4122)

$c1.Copy()
$f1 = $ws2.Range("F1")
$f1.PasteSpecial(-4122)

$d1.Copy()
$g1 = $ws2.Range("G1")
$g1.PasteSpecial(-4122)

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Drop the stray empty inline-string cell at G5
$ws2.Range("G5").ClearContents()
